# PP2 - Convert everything to containers
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openml_100")

# Rows whose Configuration (col B) moves from "snapshot" to "ensemble_snapshot",
# and which also get "Run" recorded in the PP2 column (col D).
$ensembleRowsWithRun = @(7, 16, 28, 30, 33, 48, 49, 55, 61, 62, 68, 69, 77)
foreach ($r in $ensembleRowsWithRun) {
    $ws.Cells.Item($r, 2).Value = "ensemble_snapshot"
    $ws.Cells.Item($r, 4).Value = "Run"
}

# Row that moves to "ensemble_snapshot" but keeps the PP2 column untouched.
$ws.Cells.Item(79, 2).Value = "ensemble_snapshot"

# Rows whose Configuration (col B) moves from "snapshot" to "standard" (no PP2 change).
$standardRows = @(4, 9, 20, 53, 71, 76)
foreach ($r in $standardRows) {
    $ws.Cells.Item($r, 2).Value = "standard"
}

# Row 31 moves to "standard" AND also gets "Run" recorded in the PP2 column.
$ws.Cells.Item(31, 2).Value = "standard"
$ws.Cells.Item(31, 4).Value = "Run"

# Row 34: Status (col C) moves from "progressing" to "complete", and the
# PP2 (col D) "Run" marker is cleared since the job is now finished.
$ws.Cells.Item(34, 3).Value = "complete"
$ws.Cells.Item(34, 4).ClearContents()

# Update the active selection to reflect where the user ended up working.
$ws.Range("D10").Select()

# Reposition/resize the workbook window as recorded in the saved view state.
$excel.ActiveWindow.Top = 1440
$excel.ActiveWindow.Left = 19650
$excel.ActiveWindow.Width = 18240
$excel.ActiveWindow.Height = 14400
